# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet updates ---
$rushing.Range("C2").Value = 13
$rushing.Range("D2").Value = 6
$rushing.Range("E2").Value = 17
$rushing.Range("F2").Value = 9

$rushing.Range("C3").Value = 97
$rushing.Range("D3").Value = 53
$rushing.Range("F3").Value = 34

$rushing.Range("C4").Value = 14
$rushing.Range("D4").Value = 14
$rushing.Range("F4").Value = 12

$rushing.Range("C5").Value = 12
$rushing.Range("D5").Value = 6
$rushing.Range("E5").Value = 3
$rushing.Range("F5").Value = 2

# --- Receiving sheet updates ---
$receiving.Range("C2").Value = 62
$receiving.Range("D2").Value = 49

$receiving.Range("C3").Value = 8
$receiving.Range("D3").Value = 7

$receiving.Range("C6").Value = 99
$receiving.Range("D6").Value = 70
$receiving.Range("E6").Value = 17
$receiving.Range("F6").Value = 11
$receiving.Range("G6").Value = 15
$receiving.Range("H6").Value = 10

$receiving.Range("C7").Value = 63
$receiving.Range("D7").Value = 39
$receiving.Range("E7").Value = 23
$receiving.Range("F7").Value = 11
$receiving.Range("G7").Value = 15

$receiving.Range("E8").Value = 1

$receiving.Range("C9").Value = 22
$receiving.Range("D9").Value = 12
$receiving.Range("E9").Value = 6
$receiving.Range("F9").Value = 2

$receiving.Range("C11").Value = 45
$receiving.Range("D11").Value = 28
$receiving.Range("E11").Value = 10
$receiving.Range("F11").Value = 5

$receiving.Range("C12").Value = 19
$receiving.Range("D12").Value = 15
$receiving.Range("E12").Value = 1
$receiving.Range("G12").Value = 4

# --- Active sheet / selection: Receiving becomes the active tab, selection at I12 ---
$receiving.Select()
$receiving.Range("I12").Select()
